# Generate Report for Handback
# Updates the localization-status workbook to reflect a completed handback:
#  - "Ready for handoff" status becomes "Handed back: in sync with en-US"
#    on the Overview sheet (E2/F2) and on each language sheet's Status cell (C2).
#  - Each language sheet gets its "Latest Target File" (I2) populated with the
#    source file name, as a hyperlink to the same file referenced by A2.
#  - Each language sheet gets its "Latest Handback File" (J2) populated with
#    the latest handoff xlf file name (round-tripped back).
#  - Each language sheet's "Latest Handback DateTime" (K2) is stamped.
#  - Column widths on the affected columns are widened to fit the new text.

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"
$mdFile = "868e5120-f9cc-45ae-a154-11922f54e8ed.md"
$mdHyperlinkTarget = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e14db6485b3a364582c97bbb1628d518662034ad/e2e/868e5120-f9cc-45ae-a154-11922f54e8ed.md"

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Columns.Item(5).ColumnWidth = 29.9777047293527
$wsOverview.Columns.Item(6).ColumnWidth = 29.9777047293527

# --- zh-cn sheet ---
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C2").Value = $newStatus
$wsZh.Hyperlinks.Add($wsZh.Range("I2"), $mdHyperlinkTarget, "", "", $mdFile)
$wsZh.Range("I2").Font.Underline = $true
$wsZh.Range("I2").Font.Color = 15570276
$wsZh.Range("J2").Value = "868e5120-f9cc-45ae-a154-11922f54e8ed.1349cc6be2b0898fd11a9dfe617f7cc85d069806.zh-cn.xlf"
$wsZh.Range("K2").Value = "2016-08-24 19:08:40"
$wsZh.Columns.Item(3).ColumnWidth = 29.9777047293527
$wsZh.Columns.Item(9).ColumnWidth = 40
$wsZh.Columns.Item(10).ColumnWidth = 40

# --- de-de sheet ---
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C2").Value = $newStatus
$wsDe.Hyperlinks.Add($wsDe.Range("I2"), $mdHyperlinkTarget, "", "", $mdFile)
$wsDe.Range("I2").Font.Underline = $true
$wsDe.Range("I2").Font.Color = 15570276
$wsDe.Range("J2").Value = "868e5120-f9cc-45ae-a154-11922f54e8ed.1349cc6be2b0898fd11a9dfe617f7cc85d069806.de-de.xlf"
$wsDe.Range("K2").Value = "2016-08-24 19:08:47"
$wsDe.Columns.Item(3).ColumnWidth = 29.9777047293527
$wsDe.Columns.Item(9).ColumnWidth = 40
$wsDe.Columns.Item(10).ColumnWidth = 40
